# Insert two new rows (newest-first ordering) right before the current
# row 396, pushing the existing rows 396-412 down to 398-414, then fill
# the two freshly inserted rows with the new weekly price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("396:397").Insert()

# New row 396: Plátano "Pintón"
$ws.Range("A396").Value = 5
$ws.Range("B396").Value = "Macroferia Regional de Talca"
$ws.Range("C396").Value = "Maule"
$ws.Range("D396").Value = 44509
$ws.Range("E396").Value = 7
$ws.Range("F396").Value = "Fruta"
$ws.Range("G396").Value = 100108
$ws.Range("H396").Value = "Tropicales y subtropicales"
$ws.Range("I396").Value = 100108006
$ws.Range("J396").Value = "Plátano"
$ws.Range("K396").Value = "Sin especificar"
$ws.Range("L396").Value = "Pintón"
$ws.Range("M396").Value = 1050
$ws.Range("N396").Value = 16000
$ws.Range("O396").Value = 16000
$ws.Range("P396").Value = 16000
$ws.Range("Q396").Value = "$/caja 20 kilos"
$ws.Range("R396").Value = "Ecuador"
$ws.Range("S396").Value = 800
$ws.Range("T396").Value = 20

# New row 397: Plátano "Primera Pintón"
$ws.Range("A397").Value = 5
$ws.Range("B397").Value = "Macroferia Regional de Talca"
$ws.Range("C397").Value = "Maule"
$ws.Range("D397").Value = 44509
$ws.Range("E397").Value = 7
$ws.Range("F397").Value = "Fruta"
$ws.Range("G397").Value = 100108
$ws.Range("H397").Value = "Tropicales y subtropicales"
$ws.Range("I397").Value = 100108006
$ws.Range("J397").Value = "Plátano"
$ws.Range("K397").Value = "Sin especificar"
$ws.Range("L397").Value = "Primera Pintón"
$ws.Range("M397").Value = 500
$ws.Range("N397").Value = 18000
$ws.Range("O397").Value = 18000
$ws.Range("P397").Value = 18000
$ws.Range("Q397").Value = "$/caja 20 kilos"
$ws.Range("R397").Value = "Ecuador"
$ws.Range("S397").Value = 900
$ws.Range("T397").Value = 20
